# Project Tracker: update "Phase 1" task list (Table137, rows 7-10).
#  - Task 1/2/3 placeholders replaced with the real deliverable names.
#  - Progress (% COMPLETE / DONE) values updated to reflect actual status.
#  - A brand-new row (row 10, "Final Report") is added to the phase.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "Task 1" -> "Project Plan, Team Roles and Dataset Identification"
# % COMPLETE bumped to 80% and DONE flips to 1 (literal, no longer the shared formula).
$ws.Range("A7").Value = "Project Plan, Team Roles and Dataset Identification"
$ws.Range("F7").Value = 0.8
$ws.Range("G7").Value = 1

# Row 8: "Task 2" -> "Report Draft 1"; priority downgraded from MEDIUM to LOW.
$ws.Range("A8").Value = "Report Draft 1"
$ws.Range("C8").Value = "LOW"

# Row 9: "Task 3" -> "Report Draft 2"; slips to start 43761 (END recalculates via the
# existing Table137 structured formula), now 100% complete.
$ws.Range("A9").Value = "Report Draft 2"
$ws.Range("D9").Value = 43761
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1

# Row 10: previously a blank table row, now populated with a new "Final Report" task.
$ws.Range("A10").Value = "Final Report"
$ws.Range("C10").Value = "MEDIUM"
$ws.Range("D10").Value = 43761
$ws.Range("E10").Formula = "=Table137[[#This Row],[START]]+3"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1

# Move the active selection to reflect where the author left off editing.
[void]$ws.Range("F11").Select()
